$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 320.22223
$ws.Range("J2").Value = 343
$ws.Range("L2").Value = 343
$ws.Range("N2").Value = -569
$ws.Range("H29").Value = 933.3333
$ws.Range("I29").Value = 1000
$ws.Range("K29").Value = 3000
$ws.Range("M29").Value = -2719
$ws.Range("H113").Value = 18380
$ws.Range("J113").Value = 3225
$ws.Range("L113").Value = 3225
$ws.Range("N113").Value = -9733
$ws.Range("H127").Value = 2438.4
$ws.Range("J127").Value = 2065
$ws.Range("L127").Value = 6195
$ws.Range("N127").Value = -16115
$ws.Range("H129").Value = 101543.86
$ws.Range("I129").Value = 140993.6
$ws.Range("J129").Value = 2919.5
$ws.Range("K129").Value = 422980.8
$ws.Range("L129").Value = 8758.5
$ws.Range("M129").Value = -417980.8
$ws.Range("N129").Value = -18758.5
$ws.Range("H137").Value = 11159.116
$ws.Range("I137").Value = 5693.7036
$ws.Range("J137").Value = 15630.818
$ws.Range("K137").Value = 17081.1108
$ws.Range("L137").Value = 46892.454
$ws.Range("M137").Value = -14531.1108
$ws.Range("N137").Value = -51992.454
$ws.Range("H138").Value = 11156.095
$ws.Range("I138").Value = 7887.2144
$ws.Range("K138").Value = 23661.6432
$ws.Range("M138").Value = -18521.6432

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1410738.8
$ws.Range("I32").Value = 1734365.6
$ws.Range("J32").Value = 69999.21000000001
$ws.Range("K32").Value = 1734365.6
$ws.Range("L32").Value = 69999.21000000001
$ws.Range("M32").Value = -1734078.6
$ws.Range("N32").Value = -70573.21000000001
$ws.Range("H74").Value = 9954.826999999999
$ws.Range("I74").Value = 1726.9474
$ws.Range("K74").Value = 1726.9474
$ws.Range("M74").Value = -852.9474
$ws.Range("H77").Value = 9954.826999999999
$ws.Range("I77").Value = 1726.9474
$ws.Range("K77").Value = 8634.737000000001
$ws.Range("M77").Value = -4266.737000000001
$ws.Range("H110").Value = 3845.2222
$ws.Range("I110").Value = 2172.5715
$ws.Range("K110").Value = 2172.5715
$ws.Range("M110").Value = -127.5715
$ws.Range("H132").Value = 11760.383
$ws.Range("I132").Value = 5426.5654
$ws.Range("K132").Value = 16279.6962
$ws.Range("M132").Value = -13749.6962

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1154.8
$ws.Range("I22").Value = 1218.75
$ws.Range("J22").Value = 899
$ws.Range("K22").Value = 1218.75
$ws.Range("L22").Value = 899
$ws.Range("M22").Value = -1045.75
$ws.Range("N22").Value = -1245
$ws.Range("H81").Value = 71619
$ws.Range("J81").Value = 71619
$ws.Range("L81").Value = 71619
$ws.Range("N81").Value = -73741
$ws.Range("H84").Value = 71619
$ws.Range("J84").Value = 71619
$ws.Range("L84").Value = 214857
$ws.Range("N84").Value = -225465
$ws.Range("H134").Value = 11833.441
$ws.Range("I134").Value = 5742.1665
$ws.Range("K134").Value = 17226.4995
$ws.Range("M134").Value = -14691.4995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9073.1
$ws.Range("I16").Value = 2591.5
$ws.Range("J16").Value = 34999.5
$ws.Range("K16").Value = 2591.5
$ws.Range("L16").Value = 34999.5
$ws.Range("M16").Value = -2304.5
$ws.Range("N16").Value = -35573.5
$ws.Range("H31").Value = 15541.104
$ws.Range("I31").Value = 6837.628
$ws.Range("J31").Value = 40491.066
$ws.Range("K31").Value = 6837.628
$ws.Range("L31").Value = 40491.066
$ws.Range("M31").Value = -6542.628
$ws.Range("N31").Value = -41081.066
$ws.Range("H34").Value = 15541.104
$ws.Range("I34").Value = 6837.628
$ws.Range("J34").Value = 40491.066
$ws.Range("K34").Value = 6837.628
$ws.Range("L34").Value = 40491.066
$ws.Range("M34").Value = -6635.628
$ws.Range("N34").Value = -40895.066
$ws.Range("H99").Value = 13368.533
$ws.Range("I99").Value = 4565.6665
$ws.Range("K99").Value = 4565.6665
$ws.Range("M99").Value = -3067.6665
$ws.Range("H113").Value = 9073.1
$ws.Range("I113").Value = 2591.5
$ws.Range("J113").Value = 34999.5
$ws.Range("K113").Value = 2591.5
$ws.Range("L113").Value = 34999.5
$ws.Range("M113").Value = -421.5
$ws.Range("N113").Value = -39339.5
$ws.Range("H126").Value = 13368.533
$ws.Range("I126").Value = 4565.6665
$ws.Range("K126").Value = 13696.9995
$ws.Range("M126").Value = -11226.9995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5556738
$ws.Range("I34").Value = 1398.75
$ws.Range("K34").Value = 4196.25
$ws.Range("M34").Value = -4112.25
$ws.Range("H39").Value = 2897.4
$ws.Range("J39").Value = 2000
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6588
$ws.Range("H55").Value = 2083.1667
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("H82").Value = 1000
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2594
$ws.Range("H85").Value = 1000
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1596
$ws.Range("H103").Value = 8674.375
$ws.Range("J103").Value = 13670
$ws.Range("L103").Value = 41010
$ws.Range("N103").Value = -42768
$ws.Range("H131").Value = 1463.2872
$ws.Range("J131").Value = 1484.1333
$ws.Range("L131").Value = 4452.3999
$ws.Range("N131").Value = -14532.3999
$ws.Range("H134").Value = 4015.76
$ws.Range("I134").Value = 679.2766
$ws.Range("K134").Value = 2037.8298
$ws.Range("M134").Value = 3032.1702
$ws.Range("H137").Value = 5009.6924
$ws.Range("J137").Value = 5516.7144
$ws.Range("L137").Value = 16550.1432
$ws.Range("N137").Value = -26750.1432
$ws.Range("N55").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 49999.5
$ws.Range("J28").Value = 49999.5
$ws.Range("L28").Value = 49999.5
$ws.Range("N28").Value = -50383.5
$ws.Range("H70").Value = 9798.777
$ws.Range("I70").Value = 6084.5713
$ws.Range("J70").Value = 12162.363
$ws.Range("K70").Value = 6084.5713
$ws.Range("L70").Value = 12162.363
$ws.Range("M70").Value = -5814.5713
$ws.Range("N70").Value = -12702.363
$ws.Range("H73").Value = 9798.777
$ws.Range("I73").Value = 6084.5713
$ws.Range("J73").Value = 12162.363
$ws.Range("K73").Value = 6084.5713
$ws.Range("L73").Value = 12162.363
$ws.Range("M73").Value = -5148.5713
$ws.Range("N73").Value = -14034.363
$ws.Range("H113").Value = 103018.3
$ws.Range("I113").Value = 115022.875
$ws.Range("J113").Value = 55000
$ws.Range("K113").Value = 115022.875
$ws.Range("L113").Value = 55000
$ws.Range("M113").Value = -112852.875
$ws.Range("N113").Value = -59340
$ws.Range("H122").Value = 5661.727
$ws.Range("I122").Value = 3559.5
$ws.Range("K122").Value = 10678.5
$ws.Range("M122").Value = -8228.5
$ws.Range("H132").Value = 5443.087
$ws.Range("I132").Value = 6830.143
$ws.Range("J132").Value = 3285.4443
$ws.Range("K132").Value = 20490.429
$ws.Range("L132").Value = 9856.332900000001
$ws.Range("M132").Value = -17960.429
$ws.Range("N132").Value = -14916.3329

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2871
$ws.Range("I2").Value = 2573.25
$ws.Range("J2").Value = 3665
$ws.Range("K2").Value = 2573.25
$ws.Range("L2").Value = 3665
$ws.Range("M2").Value = -2461.25
$ws.Range("N2").Value = -3889
$ws.Range("H40").Value = 8809.666999999999
$ws.Range("I40").Value = 5387.7144
$ws.Range("J40").Value = 13600.4
$ws.Range("K40").Value = 5387.7144
$ws.Range("L40").Value = 13600.4
$ws.Range("M40").Value = -5251.7144
$ws.Range("N40").Value = -13872.4
$ws.Range("H94").Value = 42499.5
$ws.Range("J94").Value = 42499.5
$ws.Range("L94").Value = 42499.5
$ws.Range("N94").Value = -43851.5
$ws.Range("H98").Value = 36000
$ws.Range("J98").Value = 36000
$ws.Range("L98").Value = 36000
$ws.Range("N98").Value = -41990
$ws.Range("H122").Value = 10599.846
$ws.Range("I122").Value = 7079.9
$ws.Range("J122").Value = 22333
$ws.Range("K122").Value = 21239.7
$ws.Range("L122").Value = 66999
$ws.Range("M122").Value = -18789.7
$ws.Range("N122").Value = -71899
$ws.Range("H132").Value = 8070.9653
$ws.Range("I132").Value = 4002.1304
$ws.Range("J132").Value = 23668.166
$ws.Range("K132").Value = 12006.3912
$ws.Range("L132").Value = 71004.49800000001
$ws.Range("M132").Value = -9476.3912
$ws.Range("N132").Value = -76064.49800000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 53321.25
$ws.Range("J56").Value = 60000
$ws.Range("L56").Value = 60000
$ws.Range("N56").Value = -61428
$ws.Range("H62").Value = 7497.5
$ws.Range("I62").Value = 9995
$ws.Range("K62").Value = 9995
$ws.Range("M62").Value = -9371
$ws.Range("H65").Value = 7497.5
$ws.Range("I65").Value = 9995
$ws.Range("K65").Value = 49975
$ws.Range("M65").Value = -46855
$ws.Range("H96").Value = 3627.4285
$ws.Range("I96").Value = 4997.5
$ws.Range("K96").Value = 4997.5
$ws.Range("M96").Value = -3624.5
$ws.Range("H107").Value = 3811.55
$ws.Range("I107").Value = 1657.2
$ws.Range("J107").Value = 5965.9
$ws.Range("K107").Value = 4971.6
$ws.Range("L107").Value = 17897.7
$ws.Range("M107").Value = -3051.6
$ws.Range("N107").Value = -21737.7
$ws.Range("H126").Value = 52070.715
$ws.Range("I126").Value = 70999
$ws.Range("J126").Value = 4750
$ws.Range("K126").Value = 212997
$ws.Range("L126").Value = 14250
$ws.Range("M126").Value = -210527
$ws.Range("N126").Value = -19190
